$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Exceptions")

# Remove the old repeated rows (4-16); clear row 1 in place (without shifting
# the remaining rows up) so the two kept rows stay addressed as A2 / A3.
$ws.Rows("4:16").Delete()
$ws.Rows("1:1").ClearContents()

# Write the updated Business-Exception messages into the two remaining rows.
$ws.Range("A2").Value = "CHANGES - SOX...eml on 01/24/2020 00:00:00 missing Server Name magic_qq_appl (Expected Server Name)"
$ws.Range("A3").Value = "CHANGES - SOX...eml on 01/24/2020 00:00:00 missing Server Name testps9023 (Expected Server Name)"

[void]$ws.Range("M6").Select()
